$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'26.446.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  +2.79%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'1.675.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  +4.14%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Formula = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Formula = "'  -0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'216.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  +4.52%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'0.5312"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  +3.00%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'  -0.09%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Formula = "'  +4.99%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Formula = "'0.06396"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'  +3.73%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'21.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  +7.59%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Formula = "'0.07792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'  +3.78%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'1.686.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'  +5.58%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Formula = "'4.506"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'  +4.07%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'0.5567"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  +3.25%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'0.0₅8352"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  +7.52%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'65.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  +3.53%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'26.470.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  +2.92%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'1.004"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  +0.07%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Formula = "'4.780"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  +4.40%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'196.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  +7.70%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'10.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  +4.19%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'6.336"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'  +5.75%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  -0.12%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'143.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  +0.07%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'0.1280"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'  +7.19%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Formula = "'7.421"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'  +1.75%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Formula = "'16.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'  +6.45%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Formula = "'1.429"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "'  +5.75%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Formula = "'0.06138"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "'  +5.55%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Formula = "'1.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  +3.50%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Formula = "'3.613"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "'  +8.41%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Formula = "'3.451"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  +4.73%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Formula = "'1.689"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  +6.81%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Formula = "'1.006"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "'  +5.14%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Formula = "'2.424"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "'  +1.85%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Formula = "'2.780"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  +3.02%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Formula = "'0.5714"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  +0.41%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Formula = "'0.01637"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "'  +4.17%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Formula = "'6.021"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Formula = "'  +7.29%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'1.069.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  +5.50%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'0.8588"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'  +3.04%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Formula = "'1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'  -0.11%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'99.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  +1.30%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Formula = "'1.822.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'  +3.64%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'0.0₈111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  +4.37%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Formula = "'57.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Formula = "'  +6.54%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Formula = "'Frax"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Formula = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Formula = "'1.004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  +0.24%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Formula = "'EnergySwap"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Formula = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Formula = "'8.119"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  +4.13%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Formula = "'0.05207"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'  +1.03%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Formula = "'Aptos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Formula = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Formula = "'6.035"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'  +5.20%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Formula = "'Mantle"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Formula = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Formula = "'0.4238"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  +0.67%  "
$ws.Range("E51").Style = "Normal"
